# Power Supply Trimming.xlsx - "Fixed AC simulation. Schematic now ready."
# Adds a voltage-reference comparison table (Price / Tolerance / Temp Coef / Noise)
# for ADR4533BRZ, LM4132CQ1MFT3.3 and ADR366BUJZ, and clears the old A1 label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "ADR4533BRZ" label that lived alone in row 1 ---
# Deleting row 1 and re-inserting a blank row 1 clears A1 while keeping every
# other row (2-15) at its original row number (so the sheet dimension becomes
# A2:I.. instead of A1:I..).
$ws.Rows(1).Delete()
$ws.Rows(1).Insert()

# --- New comparison table (rows 20-23) ---
# Values are written in the same order the original author entered them so the
# generated shared-string table lines up with the source workbook.
$ws.Range("A22").Value = "LM4132CQ1MFT3.3"
$ws.Range("A23").Value = "ADR366BUJZ"

$ws.Range("D20").Value = "Temp Coef"
$ws.Range("E20").Value = "Noise"

$ws.Range("C23").Value = "±0.13%"
$ws.Range("D23").Value = "9ppm/°C"
$ws.Range("E23").Value = "9.3µVp-p"

$ws.Range("C22").Value = "±0.2%"
$ws.Range("D22").Value = "20ppm/°C"
$ws.Range("E22").Value = "310µVp-p"

$ws.Range("B20").Value = "Price"

$ws.Range("E21").Style = "Normal"
$ws.Range("E21").Value = "2.1µVp-p"
$ws.Range("D21").Style = "Normal"
$ws.Range("D21").Value = "2ppm/°C"
$ws.Range("C21").Value = "±0.02%"

$ws.Range("A21").Value = "ADR4533BRZ"
$ws.Range("C20").Value = "Tolerance"

$ws.Range("B21").Value = 10.8
$ws.Range("B22").Value = 4.11
$ws.Range("B23").Value = 4.24

# --- Column width tweaks (author widened/rebalanced the columns for the new table) ---
$ws.Columns("A").ColumnWidth = 16.5
$ws.Columns("B").ColumnWidth = 17
$ws.Columns("C").ColumnWidth = 8.333333
$ws.Columns("D").ColumnWidth = 9
$ws.Columns("E").ColumnWidth = 8
$ws.Columns("G").ColumnWidth = 5.166667
$ws.Columns("H").ColumnWidth = 7.333333
$ws.Columns("I").ColumnWidth = 5.666667

# --- Leave the selection where the author left it ---
$ws.Range("H11").Select()
